# added ifo gdp component analysis preprocessing
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated error-table values for rows 2-10 (columns B:G), reflecting
# the refreshed ifo GDP component analysis preprocessing.
$updates = @(
    @{ Row = 2;  B = 0.09008902633495776; C = 0.2712929561332827;  D = 0.1215236909347529;  E = 0.3486024826858709;  F = 0.3459885601116726;  G = 19 }
    @{ Row = 3;  B = 0.440854564544955;   C = 0.4977539282968331;  D = 0.4961366786238803;  E = 0.7043697030848788;  F = 0.5652750901259675;  G = 18 }
    @{ Row = 4;  B = 0.6924743815712978;  C = 0.7117042034689857;  D = 0.9527927106825833;  E = 0.9761110135033737;  F = 0.7091201857914559;  G = 17 }
    @{ Row = 5;  B = 0.6938486294701675;  C = 0.696168599774547;   D = 0.7176852045027668;  E = 0.8471630330123989;  F = 0.5020058792592104;  G = 16 }
    @{ Row = 6;  B = 0.5732877799022984;  C = 0.5732877799022984;  D = 0.438377994181727;   E = 0.6621011963300829;  F = 0.3428646894649251;  G = 15 }
    @{ Row = 7;  B = 0.4864587127141264;  C = 0.4881921051911807;  D = 0.3241750264404118;  E = 0.5693637031286871;  F = 0.3070280946439436;  G = 14 }
    @{ Row = 8;  B = 0.4121643343351353;  C = 0.4244592671431573;  D = 0.2374407211537076;  E = 0.4872788946319219;  F = 0.2705390721203952;  G = 13 }
    @{ Row = 9;  B = 0.4514641206974861;  C = 0.4514641206974861;  D = 0.2461149259357549;  E = 0.4960997137025528;  F = 0.2181328903808997;  G = 9 }
    @{ Row = 10; B = 0.3763340750174313;  C = 0.3763340750174313;  D = 0.1542216358128203;  E = 0.3927106260503022;  F = 0.1254706130613591;  G = 5 }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Range("B$r").Value = $u.B
    $ws.Range("C$r").Value = $u.C
    $ws.Range("D$r").Value = $u.D
    $ws.Range("E$r").Value = $u.E
    $ws.Range("F$r").Value = $u.F
    $ws.Range("G$r").Value = $u.G
}
